# Lattice-multiplication worksheet refresh: every exercise cell in the
# single 5x3 table gets a new "A x B" problem (and its matching lattice
# scaffold digits) per the "Update master to output generated at
# 503736d" regeneration. Each cell holds ONE run with five <w:t> lines
# separated by <w:br/>: the "A x B" title, the spaced-out digits of B,
# the "----" divider, and finally the two digits of A (each paired with
# a lattice "|    |" box). Because several of those scaffold lines repeat
# verbatim across different cells (e.g. "  6    7" / "4|    |"), a plain
# document-wide Find/Replace would be ambiguous; instead we target each
# table cell's own Range directly and rebuild its run via InsertXML so
# the five-line, <w:br/>-separated structure (and the xml:space=
# "preserve" flag on the space-padded lines) is reproduced exactly.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row/Col (1-based) -> new ("A x B", "  b1    b2", "a1|    |", "a2|    |")
$newProblems = @(
    @{ Row = 1; Col = 1; Title = "14 x 84"; Digits = "  8    4"; Left1 = "1|    |"; Left2 = "4|    |" },
    @{ Row = 1; Col = 2; Title = "95 x 36"; Digits = "  3    6"; Left1 = "9|    |"; Left2 = "5|    |" },
    @{ Row = 1; Col = 3; Title = "23 x 82"; Digits = "  8    2"; Left1 = "2|    |"; Left2 = "3|    |" },
    @{ Row = 2; Col = 1; Title = "69 x 46"; Digits = "  4    6"; Left1 = "6|    |"; Left2 = "9|    |" },
    @{ Row = 2; Col = 2; Title = "66 x 77"; Digits = "  7    7"; Left1 = "6|    |"; Left2 = "6|    |" },
    @{ Row = 2; Col = 3; Title = "92 x 71"; Digits = "  7    1"; Left1 = "9|    |"; Left2 = "2|    |" },
    @{ Row = 3; Col = 1; Title = "22 x 90"; Digits = "  9    0"; Left1 = "2|    |"; Left2 = "2|    |" },
    @{ Row = 3; Col = 2; Title = "74 x 92"; Digits = "  9    2"; Left1 = "7|    |"; Left2 = "4|    |" },
    @{ Row = 3; Col = 3; Title = "67 x 87"; Digits = "  8    7"; Left1 = "6|    |"; Left2 = "7|    |" },
    @{ Row = 4; Col = 1; Title = "63 x 43"; Digits = "  4    3"; Left1 = "6|    |"; Left2 = "3|    |" },
    @{ Row = 4; Col = 2; Title = "56 x 45"; Digits = "  4    5"; Left1 = "5|    |"; Left2 = "6|    |" },
    @{ Row = 4; Col = 3; Title = "30 x 99"; Digits = "  9    9"; Left1 = "3|    |"; Left2 = "0|    |" },
    @{ Row = 5; Col = 1; Title = "85 x 17"; Digits = "  1    7"; Left1 = "8|    |"; Left2 = "5|    |" },
    @{ Row = 5; Col = 2; Title = "35 x 29"; Digits = "  2    9"; Left1 = "3|    |"; Left2 = "5|    |" },
    @{ Row = 5; Col = 3; Title = "17 x 18"; Digits = "  1    8"; Left1 = "1|    |"; Left2 = "7|    |" }
)

foreach ($p in $newProblems) {
    $innerDoc = "<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" +
                "<w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr>" +
                "<w:t>$($p.Title)</w:t><w:br/>" +
                "<w:t xml:space=`"preserve`">$($p.Digits)</w:t><w:br/>" +
                "<w:t xml:space=`"preserve`">  ----</w:t><w:br/>" +
                "<w:t>$($p.Left1)</w:t><w:br/>" +
                "<w:t>$($p.Left2)</w:t>" +
                "</w:r></w:p></w:body></w:document>"

    $pkgXml = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
              "<pkg:part pkg:name=`"/word/document.xml`" " +
              "pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
              "<pkg:xmlData>$innerDoc</pkg:xmlData></pkg:part></pkg:package>"

    $cellRange = $t.Cell($p.Row, $p.Col).Range
    $cellRange.InsertXML($pkgXml)
}
